$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = -8.200899999999997
$ws.Range("D18").Value = -8.848999999999995
$ws.Range("D20").Value = -8.264499999999998
$ws.Range("D27").Value = -7.683300000000002
$ws.Range("D69").Value = -7.546799999999998
$ws.Range("D76").Value = -7.692299999999999
$ws.Range("D82").Value = -8.45240000000001
